$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------
# 1. Drop the stray "_GoBack" bookmark that used to sit in front of the
#    "Реферат" heading paragraph.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2. "Haswell" -> "Ivy Bridge" (processor family table, row 1) and mark
#    the new run/paragraph text as English (en-US).
# ---------------------------------------------------------------------
$specs = $d.Tables.Item(2)
$familyPara = $specs.Cell(1, 2).Range.Paragraphs.Item(1).Range
$familyXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="' + $wNs + '"><w:body>' +
    '<w:p w:rsidR="00E75343" w:rsidRDefault="00224B4B">' +
    '<w:pPr><w:ind w:left="140"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Ivy Bridge</w:t></w:r>' +
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$familyPara.InsertXML($familyXml)

# ---------------------------------------------------------------------
# 3. "4" -> "8" (number-of-cores table row) and re-add a fresh "_GoBack"
#    bookmark right after the new run, at the end of the paragraph.
# ---------------------------------------------------------------------
$coresPara = $specs.Cell(2, 2).Range.Paragraphs.Item(1).Range
$coresXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="' + $wNs + '"><w:body>' +
    '<w:p w:rsidR="00E75343" w:rsidRDefault="00224B4B">' +
    '<w:pPr><w:ind w:left="140"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>8</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$coresPara.InsertXML($coresXml)
